$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column letters being updated, in order: B, C, E, F, G, I, L, O
$cols = @(2, 3, 5, 6, 7, 9, 12, 15)

$data = @{
    2 = @(1.359058946454525, 0.2396082070310683, 0.08610189442810423, 0.4443680307746121, 0.002418976509366153, 0.5353633724865219, 0.2268197972997825, 2.027512517376096)
    3 = @(1.216404908066409, 0.2219112585119376, 0.08718609321774728, 0.3878228170618172, 0.002421721691165057, 0.5479572496751448, 0.216317066434911, 2.062135403644874)
    4 = @(1.128653607467186, 0.2109989776520251, 0.08790745916207143, 0.3531389305169483, 0.002423495970071075, 0.5561889617213964, 0.2099696181120549, 2.08541484571542)
    5 = @(1.092856019552016, 0.2065408226471845, 0.08821541725921644, 0.3390132514313251, 0.002424241380730489, 0.5596686720773292, 0.2074085404963597, 2.095408260974494)
    6 = @(1.086909617119602, 0.2057998758603219, 0.08826739872519784, 0.336668177824194, 0.002424366509108832, 0.5602540306788164, 0.2069848220051966, 2.09709823434801)
    7 = @(1.128170980363848, 0.2109388987798297, 0.0879115557267669, 0.3529483938344953, 0.002423505932244739, 0.5562353836463352, 0.2099349749101549, 2.085547569915491)
    8 = @(1.309906609015286, 0.2335160945901862, 0.08646417241040361, 0.4248636149813478, 0.002419904672989015, 0.5396020137661779, 0.2231774891500322, 2.039030279056448)
    9 = @(1.664928018489206, 0.2774105350216587, 0.08406758491183375, 0.5661985755041457, 0.00241354357488544, 0.5109547794275127, 0.249947022753517, 1.963902840553331)
    10 = @(1.924847240330166, 0.3094154223363716, 0.08257618378541487, 0.6702781546542269, 0.002409293120915472, 0.4923430926032832, 0.2701018426284065, 1.918597152611611)
    11 = @(2.042876517011848, 0.3239196902138985, 0.08195621710470213, 0.7176906081379002, 0.002407450427421083, 0.4844081214156546, 0.2793765254820073, 1.90015112886131)
    12 = @(2.087539152566308, 0.3294039066941821, 0.08172986346830236, 0.7356546913071611, 0.002406765644078448, 0.4814800733261748, 0.282903815182209, 1.893478662049588)
    13 = @(2.077921737755219, 0.3282231541472527, 0.08177823846635235, 0.7317853510981394, 0.002406912547089245, 0.4821072622681264, 0.2821434771315552, 1.8949017685499)
    14 = @(2.046551603697878, 0.3243710469701284, 0.08193742619118005, 0.7191683204515869, 0.002407393829699434, 0.4841656890954003, 0.2796664142757379, 1.899595906635867)
    15 = @(2.02733218897805, 0.3220104376885047, 0.08203602920709052, 0.7114413442032514, 0.002407690321016849, 0.4854365397899194, 0.2781511143079598, 1.902511961299339)
    16 = @(1.917129332885111, 0.3084664002792863, 0.08261787725029812, 0.6671810134426437, 0.002409415368364097, 0.492872385107761, 0.2694978481400483, 1.919846312015977)
    17 = @(1.849468135646646, 0.3001432649319042, 0.08298980371629838, 0.6400460337125793, 0.002410496856968758, 0.4975704081238401, 0.2642164695275255, 1.931035631526001)
    18 = @(1.810531647619996, 0.2953508686472333, 0.08320923002115777, 0.6244449056556647, 0.002411127456360584, 0.5003226026930196, 0.2611887576806566, 1.937675006418687)
    19 = @(1.797345136284434, 0.2937273716351285, 0.08328446916646826, 0.619163680173358, 0.002411342437706798, 0.5012630292417999, 0.2601653464878808, 1.939957903111349)
    20 = @(1.856672827104205, 0.3010298122596282, 0.08294964182494802, 0.642933953830422, 0.002410380845580079, 0.4970651168171756, 0.2647776472997663, 1.929823429335798)
    21 = @(2.055766673227538, 0.3255027303941631, 0.08189044049617245, 0.7228739723491628, 0.002407252112848638, 0.483558993684154, 0.2803935767734771, 1.898208627576594)
    22 = @(2.185695540762879, 0.3414490431171373, 0.08124723842313308, 0.7751780083420101, 0.002405283078438344, 0.4751794834043128, 0.290687859444489, 1.879369535353533)
    23 = @(2.116368271896647, 0.3329427127422377, 0.08158603799594211, 0.7472568307830727, 0.002406327075819256, 0.4796107341222609, 0.2851855527812148, 1.88925701490507)
    24 = @(1.853415701679694, 0.3006290269741783, 0.08296778157041729, 0.6416283278902171, 0.002410433266624475, 0.4972933994337918, 0.2645239119636926, 1.930370822862386)
    25 = @(1.569039970460381, 0.2655779075772671, 0.0846686266161143, 0.5279251897347166, 0.002415189821875918, 0.5182779927988417, 0.2426195155729971, 1.98249568983141)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($r, $cols[$i]).Value = $vals[$i]
    }
}
